$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UserList")
$ws2 = $wb.Worksheets.Item("RequestList")

# --- Column G (Employee Status): switch from numeric 1/0 to text "t"/"f" ---
$ws1.Range("G2").Value = "t"
$ws1.Range("G2").Style = "Normal"
$ws1.Range("G3").Value = "t"
$ws1.Range("G4").Value = "t"
$ws1.Range("G5").Value = "t"
$ws1.Range("G6").Value = "f"
$ws1.Range("G7").Value = "f"
$ws1.Range("G8").Value = "f"
$ws1.Range("G9").Value = "f"

# --- Column E (Card #): apply integer number format to the whole column ---
$ws1.Range("H2").NumberFormat = "0.00E+00"
$ws1.Range("E2:E9").NumberFormat = "0"

# --- Column widths: split the D:E combined width into individual widths ---
$ws1.Columns("D").ColumnWidth = 10
$ws1.Columns("E").ColumnWidth = 12.5703125

# --- Selection / active sheet ---
$ws2.Range("F3").Select()

$ws1.Range("F6").Select()
$ws1.Activate()
